# Week 15 logged / Week 16 simulated - update per-play logs and season totals

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet: append newly logged play-by-play yardage values
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$c = $ws.Cells.Item(2,2)
$c.Value = $c.Text + " 2 -2 1 3 19 7 6 3 3 6 3 5 8 6 3 3 0 1 6 2 9 4 1 -3"

$c = $ws.Cells.Item(2,3)
$c.Value = $c.Text + " 4 -7 -2 8 1 4 4 -1 2 4 24 2 4 1 1 6 4 3 12 2 6 2 5 0"

$c = $ws.Cells.Item(3,2)
$c.Value = $c.Text + " -1 12 14 11 17 12 7 5 5 10 5 13 8 6 16 14 16 12 6 7 3 12 9 12 15"

$c = $ws.Cells.Item(3,3)
$c.Value = $c.Text + " 5 10 6 6 2 4 3 4 4 9 13 6 6 13 5 25 11 3 6"

# ---------------------------------------------------------------
# ST sheet: append newly logged return-yardage values
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$c = $ws.Cells.Item(4,2)
$c.Value = $c.Text + " 62"

$c = $ws.Cells.Item(5,2)
$c.Value = $c.Text + " 24"

$c = $ws.Cells.Item(6,2)
$c.Value = $c.Text + " 14 23"

$c = $ws.Cells.Item(3,4)
$c.Value = $c.Text + " 42 49 41 45"

$c = $ws.Cells.Item(4,4)
$c.Value = $c.Text + " 20 0 0 0"

$c = $ws.Cells.Item(5,4)
$c.Value = $c.Text + " 29 0 0 0 2"

# ---------------------------------------------------------------
# OFF sheet: updated season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")
$ws.Cells.Item(2,2).Value = 5
$ws.Cells.Item(2,3).Value = 351
$ws.Cells.Item(2,4).Value = 33
$ws.Cells.Item(2,6).Value = 85
$ws.Cells.Item(2,7).Value = 109
$ws.Cells.Item(2,8).Value = 7
$ws.Cells.Item(2,10).Value = 48
$ws.Cells.Item(2,14).Value = 31
$ws.Cells.Item(2,15).Value = 36
$ws.Cells.Item(2,16).Value = 24

$ws.Cells.Item(3,2).Value = 15
$ws.Cells.Item(3,3).Value = 343
$ws.Cells.Item(3,5).Value = 52
$ws.Cells.Item(3,6).Value = 218
$ws.Cells.Item(3,7).Value = 81
$ws.Cells.Item(3,8).Value = 36
$ws.Cells.Item(3,9).Value = 99
$ws.Cells.Item(3,10).Value = 103
$ws.Cells.Item(3,12).Value = 505
$ws.Cells.Item(3,13).Value = 337
$ws.Cells.Item(3,17).Value = 909

# ---------------------------------------------------------------
# DEF sheet: updated season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")
$ws.Cells.Item(2,3).Value = 378
$ws.Cells.Item(2,5).Value = 15
$ws.Cells.Item(2,6).Value = 120
$ws.Cells.Item(2,8).Value = 8
$ws.Cells.Item(2,10).Value = 54
$ws.Cells.Item(2,15).Value = 37
$ws.Cells.Item(2,16).Value = 27

$ws.Cells.Item(3,2).Value = 19
$ws.Cells.Item(3,3).Value = 311
$ws.Cells.Item(3,6).Value = 193
$ws.Cells.Item(3,7).Value = 63
$ws.Cells.Item(3,8).Value = 53
$ws.Cells.Item(3,9).Value = 115
$ws.Cells.Item(3,10).Value = 98
$ws.Cells.Item(3,12).Value = 529
$ws.Cells.Item(3,13).Value = 338
$ws.Cells.Item(3,17).Value = 974

# ---------------------------------------------------------------
# ST sheet: updated season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")
$ws.Cells.Item(2,2).Value = 160
$ws.Cells.Item(2,4).Value = 97
$ws.Cells.Item(2,6).Value = 134
$ws.Cells.Item(2,7).Value = 127
$ws.Cells.Item(2,10).Value = 66
$ws.Cells.Item(2,11).Value = 63
$ws.Cells.Item(2,12).Value = 35
$ws.Cells.Item(2,13).Value = 25
$ws.Cells.Item(3,2).Value = 101

# ---------------------------------------------------------------
# TURNS sheet: updated season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")
$ws.Cells.Item(3,2).Value = 9
$ws.Cells.Item(3,4).Value = 12
$ws.Cells.Item(3,5).Value = 13

# ---------------------------------------------------------------
# PEN sheet: updated season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")
$ws.Cells.Item(3,2).Value = 30
$ws.Cells.Item(4,4).Value = 21
